# The BMKG "Juni 2014" workbook originally has a single sheet
# ("Data Harian - Table") holding station metadata plus a daily-readings
# table in A9:K39. The edit duplicates that table (header + 30 daily rows)
# as plain values/formats into a brand-new worksheet ("Sheet1") placed
# right after it, then leaves that new sheet active/selected - i.e. a
# "new dataset preprocessing" sheet built from a copy of the daily table.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(1)

# Remember the table that gets copied.
$table = $src.Range("A9:K39")

# Insert the new worksheet right after the source sheet; Excel names the
# first fresh sheet it mints "Sheet1", matching sheetId 2 / rId2 in the
# saved workbook.
$dst = $wb.Worksheets.Add($null, $src)

# Copy values + formatting (border/alignment styles s=1 header, s=2 body)
# from the source table straight into the new sheet starting at A1.
$table.Copy($dst.Range("A1"))

# Match the row height Excel computed for the wrapped/bordered body rows
# once they landed in the new sheet.
$dst.Range("A2:K31").RowHeight = 28.8

# Restore the source sheet's own selection (it was left with the table
# selected after being copied from) ...
[void]$src.Range("A9:K39").Select()

# ... then switch to and select the whole of the freshly pasted table on
# the new sheet, which is what ends up as the active tab/selection.
$dst.Activate()
[void]$dst.Range("A1:K31").Select()
